$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: average of column J (k value) across the 10 data rows ---
$ws.Range("A12:Z12").Value = $null
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- Row 14: Average of SW(S*)/SW(OPT) ---
$ws.Range("A14:Z14").Value = $null
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108
$ws.Rows(14).RowHeight = 15.6

# --- Row 15: Average of SC(S*)/SC(OPT) ---
$ws.Range("A15:Z15").Value = $null
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Rows(15).RowHeight = 15.6

# --- Row 16: Worst of SW(S*)/SW(OPT) ---
$ws.Range("A16:Z16").Value = $null
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Rows(16).RowHeight = 15.6

# --- Row 17: Worst of SC(S*)/SC(OPT) ---
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$ws.Rows(17).RowHeight = 15.6

# Copy B14's format (bold, 12pt, vertically centered) onto B15:B17 so that
# they all share the exact same cell style without generating extra
# intermediate style-table entries.
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection matching final workbook state ---
$ws.Range("A14:B17").Select()
